# Updates cryptos list figures (price/volume) per the Dec 18 2023 GitHub Actions refresh.
# Two coin pairs also swap ranking position (rows 13/14 and 34/35).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '41.608.45'
$ws.Range('E2').Value = '  -1.52%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '2.168.88'
$ws.Range('E3').Value = '  -2.93%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  +0.00%  '

# Row 5: BNB
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.80'
$ws.Range('E5').Value = '  -2.29%  '

# Row 6: XRP
$ws.Range('E6').Value = '  -2.94%  '

# Row 7: Solana
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '72.15'
$ws.Range('E7').Value = '  -3.14%  '

# Row 8: USDC
$ws.Range('E8').Value = '  -0.07%  '

# Row 9: Cardano
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.578'
$ws.Range('E9').Value = '  -4.59%  '

# Row 10: Avalanche
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.66'
$ws.Range('E10').Value = '  -7.14%  '

# Row 11: Dogecoin
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0904'
$ws.Range('E11').Value = '  -6.12%  '

# Row 12: OKB
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.36'
$ws.Range('E12').Value = '  -3.95%  '

# Row 13: Polkadot -> TRON
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0998'
$ws.Range('E13').Value = '  -3.73%  '

# Row 14: TRON -> Polkadot
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.68'
$ws.Range('E14').Value = '  -4.30%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range('D15').Value = '2.491.02'
$ws.Range('E15').Value = '  -3.07%  '

# Row 16: Chainlink
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.25'
$ws.Range('E16').Value = '  -0.68%  '

# Row 17: WrappedEther
$ws.Range('D17').Value = '2.151.84'
$ws.Range('E17').Value = '  -4.43%  '

# Row 18: Polygon
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.776'
$ws.Range('E18').Value = '  -7.59%  '

# Row 19: WrappedBTC
$ws.Range('D19').Value = '41.488.64'
$ws.Range('E19').Value = '  -1.48%  '

# Row 20: ShibaInu
$ws.Range('E20').Value = '  -3.40%  '

# Row 21: Litecoin
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '69.97'
$ws.Range('E21').Value = '  -4.08%  '

# Row 22: Uniswap
$ws.Range('E22').Value = '  -7.34%  '

# Row 23: InternetComputer(DFINITY)
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.95'
$ws.Range('E23').Value = '  -10.72%  '

# Row 24: BitcoinCash
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '226.08'
$ws.Range('E24').Value = '  -2.09%  '

# Row 25: ImmutableX
$ws.Range('E25').Value = '  -4.03%  '

# Row 26: Dai
$ws.Range('E26').Value = '  -0.04%  '

# Row 27: Cosmos
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.70'
$ws.Range('E27').Value = '  -6.13%  '

# Row 29: PancakeSwap
$ws.Range('E29').Value = '  -4.23%  '

# Row 30: Toncoin
$ws.Range('E30').Value = '  -1.59%  '

# Row 31: Monero
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '170.86'
$ws.Range('E31').Value = '  +2.27%  '

# Row 32: EthereumClassic
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '19.80'
$ws.Range('E32').Value = '  -4.18%  '

# Row 33: InjectiveProtocol
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '32.52'
$ws.Range('E33').Value = '  +8.03%  '

# Row 34: Filecoin -> Hedera
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0772'
$ws.Range('E34').Value = '  -4.05%  '

# Row 35: Hedera -> Filecoin
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.38'
$ws.Range('E35').Value = '  -4.82%  '

# Row 36: Stellar
$ws.Range('E36').Value = '  -3.78%  '

# Row 37: RenderToken
$ws.Range('E37').Value = '  -1.09%  '

# Row 38: Kaspa
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.102'
$ws.Range('E38').Value = '  -7.41%  '

# Row 39: VeChain
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0304'
$ws.Range('E39').Value = '  -0.51%  '

# Row 40: Celestia
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '12.11'
$ws.Range('E40').Value = '  -8.66%  '

# Row 41: LidoDAOToken
$ws.Range('E41').Value = '  -2.14%  '

# Row 42: THORChain
$ws.Range('E42').Value = '  -6.22%  '

# Row 43: MultiversX
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '58.69'
$ws.Range('E43').Value = '  -9.95%  '

# Row 44: FraxShare
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.42'
$ws.Range('E44').Value = '  -3.61%  '

# Row 45: Algorand
$ws.Range('E45').Value = '  -5.78%  '

# Row 46: Cronos
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0964'
$ws.Range('E46').Value = '  -3.91%  '

# Row 47: Aave
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '97.38'
$ws.Range('E47').Value = '  -7.04%  '

# Row 48: ARBITRUM
$ws.Range('E48').Value = '  -4.22%  '

# Row 49: TrustWalletToken
$ws.Range('E49').Value = '  -5.28%  '

# Row 50: NEARProtocol
$ws.Range('E50').Value = '  -7.74%  '

# Row 51: HuobiToken
$ws.Range('E51').Value = '  -2.33%  '
